$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two main input values (Q18, Q19) that drive the recalculation
$ws.Range("Q18").Value = 90
$ws.Range("Q19").Value = 75

# New "A" label and formula in P21 / Q21
$ws.Range("P21").Value = "A"
$ws.Range("Q21").Formula = "=Q20*(Q18+Q19)"

# Clear out old P22/Q22 content (d2 label + 2.5 value)
$ws.Range("P22").ClearContents()
$ws.Range("Q22").Clear()

# Update R27 formula to use new Q21 cell
$ws.Range("R27").Formula = "=ROUND(C2/Q21,2)"

# Update safety-factor input values
$ws.Range("U25").Value = 160
$ws.Range("U26").Value = 190

# Add new safety factor formulas
$ws.Range("V28").Formula = "=ROUND(U26/U28,2)"
$ws.Range("V29").Formula = "=ROUND(U25/U29,2)"

# Update the selection to match the saved view state
$ws.Range("X8").Select()
